$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename a couple of existing NAME_EU values (Euskera translations)
$ws.Range("C7").Value = "Sofritoa"

# Fill in the new ARTICULO_EU column (F) for all data rows.
# Most rows use "before" (corresponds to ARTICULO_EN="before"/ARTICULO_ES="de"),
# a couple of "con" rows use "rekin", and row 24 uses "after".
$ws.Range("F2").Value = "before"
$ws.Range("F3").Value = "before"
$ws.Range("F4").Value = "before"
$ws.Range("F5").Value = "before"
$ws.Range("F6").Value = "before"
$ws.Range("F7").Value = "before"
$ws.Range("F8").Value = "before"
$ws.Range("F9").Value = "before"
$ws.Range("F10").Value = "rekin"
$ws.Range("F11").Value = "before"
$ws.Range("F12").Value = "before"
$ws.Range("F13").Value = "before"
$ws.Range("F14").Value = "before"
$ws.Range("F15").Value = "rekin"
$ws.Range("F16").Value = "before"
$ws.Range("F17").Value = "before"
$ws.Range("F18").Value = "before"
$ws.Range("F19").Value = "before"
$ws.Range("F20").Value = "before"
$ws.Range("F21").Value = "before"
$ws.Range("F22").Value = "before"
$ws.Range("F23").Value = "before"
$ws.Range("F24").Value = "after"
$ws.Range("F25").Value = "before"
$ws.Range("F26").Value = "before"
$ws.Range("F27").Value = "before"
$ws.Range("F28").Value = "before"
$ws.Range("F29").Value = "before"
$ws.Range("F30").Value = "before"

# Rename the second Euskera translation last so the new shared-string order
# matches (Sofritoa, rekin, Betelana)
$ws.Range("C23").Value = "Betelana"

# Update the active selection to match the author's last cursor position
$ws.Range("I12").Select() | Out-Null
